# Add one more day/place row to the mapping sheet:
#   A4 = A3 + 1  (next day's date, inherits A3's date style)
#   B4 = "Murdeshwar" (same place name already used in B2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4, column A: formula referencing the previous day
$ws.Range("A4").Formula = "=A3+1"

# Copy A3's formatting (date number format) onto A4 so it keeps style index 1
# instead of Excel minting a brand-new number format/style.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New row 4, column B: place name (reuses existing shared string "Murdeshwar")
$ws.Range("B4").Value = "Murdeshwar"

# Excel leaves the active selection on the next empty row after data entry
$ws.Range("A5").Select()
